$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: I1 = "I0", J1 = "IF" (same header style as existing headers, copied from H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-11: I = 1 (except row 11 = 3), J = copy of H (except row 11 = 4)
$data = @(
    @(2, 1, 6),
    @(3, 1, 5),
    @(4, 1, 6),
    @(5, 1, 3),
    @(6, 1, 6),
    @(7, 1, 5),
    @(8, 1, 5),
    @(9, 1, 6),
    @(10, 1, 2),
    @(11, 3, 4)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}
